$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 64; this shifts existing rows 64-134 down to 65-135
# and expands the used range to A1:R135, matching the target diff.
$ws.Rows(64).Insert()

# Populate the newly inserted row 64 with a new weekly price observation.
$ws.Cells.Item(64, 1).Value = 8
$ws.Cells.Item(64, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(64, 3).Value = "Coquimbo"
$ws.Cells.Item(64, 4).Value = 44880
$ws.Cells.Item(64, 5).Value = 4
$ws.Cells.Item(64, 6).Value = 100112052
$ws.Cells.Item(64, 7).Value = "Albahaca"
$ws.Cells.Item(64, 8).Value = "Sin especificar"
$ws.Cells.Item(64, 9).Value = "Primera"
$ws.Cells.Item(64, 10).Value = 1100
$ws.Cells.Item(64, 11).Value = 3800
$ws.Cells.Item(64, 12).Value = 4000
$ws.Cells.Item(64, 13).Value = 3900
$ws.Cells.Item(64, 14).Value = "`$/paquete"
$ws.Cells.Item(64, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(64, 16).Value = 3900
$ws.Cells.Item(64, 17).Value = 1
$ws.Cells.Item(64, 18).Value = "Hortaliza"
